# The workbook rows were re-sorted (by Id, column A) after a few GPS / date
# corrections, which causes a 3-way cyclic rotation of the data held in
# rows 19, 21 and 22 (row numbers/formatting stay put, only the record
# contents move): old row19 -> row21, old row21 -> row22, old row22 -> row19.
#
# Below we just write each destination row's final values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 19 (becomes what used to be row 22's record) ----
$ws.Range("A19").Value = 111559701
$ws.Range("B19").Value = 12249
$ws.Range("D19").Value = "EN"
$ws.Range("E19").Value = 101283
$ws.Range("F19").Value = "Djupsvart brunbagge"
$ws.Range("G19").Value = "Melandrya dubia"
$ws.Range("H19").Value = "(Schaller, 1783)"
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("Q19").Value = 523950.9321204902
$ws.Range("R19").Value = 6934675.944620069
$ws.Range("Y19").Value = "'2023-08-17"
$ws.Range("AA19").Value = "'2023-08-17"
$ws.Range("AC19").Value = "Kläckhål med svartfärgade larvgångar på björkhögstubbe med levande fnösktickor. Naturskog norr om Vattensjöarna"

# ---- Row 21 (becomes what used to be row 19's record) ----
$ws.Range("A21").Value = 111560043
$ws.Range("Q21").Value = 523949.236686704
$ws.Range("R21").Value = 6934654.704083432

# ---- Row 22 (becomes what used to be row 21's record) ----
$ws.Range("A22").Value = 111560058
$ws.Range("B22").Value = 78578
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6458
$ws.Range("F22").Value = "Lunglav"
$ws.Range("G22").Value = "Lobaria pulmonaria"
$ws.Range("H22").Value = "(L.) Hoffm."
$ws.Range("L22").Value = ""
$ws.Range("M22").Value = ""
$ws.Range("Q22").Value = 523906.9737172622
$ws.Range("R22").Value = 6934619.326478666
$ws.Range("Y22").Value = "'2023-08-18"
$ws.Range("AA22").Value = "'2023-08-18"
$ws.Range("AC22").Value = ""
